$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Not Started"
$ws.Range("B3").Value = "Go Live"
$ws.Range("B4").Value = "Go Live"
$ws.Range("B5").Value = "Go Live"
$ws.Range("B6").Value = "Go Live"
$ws.Range("B7").Value = "Go Live"
$ws.Range("B8").Value = "Go Live"
$ws.Range("B12").Value = "not_started"
$ws.Range("B13").Value = "not_started"
$ws.Range("B14").Value = "not_started"
$ws.Range("B15").Value = "not_started"
$ws.Range("B16").Value = "not_started"
$ws.Range("B18").Value = "not_started"
$ws.Range("B19").Value = "not_started"
$ws.Range("B28").Value = "Go Live"
$ws.Range("B29").Value = "not_started"
$ws.Range("B31").Value = "not_started"
$ws.Range("B32").Value = "not_started"
$ws.Range("B47").Value = "not_started"
$ws.Range("B49").Value = "not_started"
$ws.Range("B50").Value = "not_started"
$ws.Range("B51").Value = "not_started"
